# Adição da função de gerar etiquetas e exibir
# Renomeia a planilha existente e adiciona a nova planilha "Maggiore Modas"
# com os mesmos cabeçalhos, pronta para receber os dados da API do cliente.

$wb = $excel.ActiveWorkbook

# Renomeia a primeira planilha
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Brilha Natal"

# Adiciona a nova planilha logo depois da primeira
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Maggiore Modas"

# Copia o cabeçalho (com formatação) da primeira planilha para a nova
$ws1.Range("A1:F1").Copy($ws2.Range("A1:F1"))

# Ajusta as margens de página da nova planilha (padrão métrico)
$ws2.PageSetup.LeftMargin = 36.850393728
$ws2.PageSetup.RightMargin = 36.850393728
$ws2.PageSetup.TopMargin = 56.692913399999995
$ws2.PageSetup.BottomMargin = 56.692913399999995
$ws2.PageSetup.HeaderMargin = 22.67716464
$ws2.PageSetup.FooterMargin = 22.67716464

# Restaura a seleção da primeira planilha (linha inteira selecionada)
$ws1.Rows.Item(1).Select() | Out-Null

# Ativa a nova planilha e posiciona a seleção em I5
$ws2.Activate()
$ws2.Range("I5").Select() | Out-Null
